$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two extra sample rows (R2 / Switch1), keep only one data row ---
$ws.Range("A3:A4").EntireRow.Delete()

# --- Insert two new leading columns for "modelo" and "serie" ---
$ws.Range("A1:B1").EntireColumn.Insert()

# Copy the existing header formatting (bold, centered, bordered) onto the new header cells
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# --- Header row ---
$ws.Range("A1").Value = "modelo"
$ws.Range("B1").Value = "serie"
$ws.Range("C1").Value = "puerto"
$ws.Range("D1").Value = "baudios"
$ws.Range("E1").Value = "nombre"
$ws.Range("F1").Value = "usuario"
$ws.Range("G1").Value = "contrasena"
$ws.Range("H1").Value = "dominio"

# --- Data row (the device that actually connects now) ---
$ws.Range("A2").Value = "CISCO2901/K9"
$ws.Range("B2").Value = "FTX153782SQ"
$ws.Range("C2").Value = "COM8"
$ws.Range("D2").Value = 9600
$ws.Range("E2").Value = "RouterPrueba"
$ws.Range("F2").Value = "admin"
$ws.Range("G2").Value = "cisco123"
$ws.Range("H2").Value = "lab.local"

# --- Stray formatted cell left from editing, matches last selection ---
$ws.Range("H7").Font.Underline = 1
$ws.Range("H7").Select()
